$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.975.78"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.748.58"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'234.85"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'0.9979"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.5176"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D8").Value = "'0.2822"
$ws.Range("E8").Value = "  +8.16%  "
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "'0.06137"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "1.744.19"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'0.07020"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "'15.50"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "'0.6448"
$ws.Range("E14").Value = "  +6.26%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'77.14"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "'0.9979"
$ws.Range("D18").Value = "'0.9982"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "25.984.05"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'11.52"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "'0.000006633"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "1.965.98"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'4.144"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("D24").Value = "'8.588"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "'5.162"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'140.48"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "'1.499"
$ws.Range("E27").Value = "  +2.76%  "
$ws.Range("D28").Value = "'1.844"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'103.25"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'0.08316"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").Value = "'3.653"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").Value = "'3.434"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").Value = "'0.04427"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "'2.606"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "'0.9896"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "'0.6125"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "'2.687"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'0.01577"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'1.943"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "'0.9971"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "'100.79"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("D43").Value = "'0.3877"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").Value = "'0.7344"
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").Value = "'5.000"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").Value = "'0.05461"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "'6.378"
$ws.Range("E47").Value = "  +7.58%  "
$ws.Range("D48").Value = "'0.1122"
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("D49").Value = "'52.78"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").Value = "'29.96"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Value = "'0.3431"
$ws.Range("E51").Value = "  +0.16%  "
